$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 457 (existing rows 457..485 shift down to 459..487)
$ws.Rows.Item(457).Insert()
$ws.Rows.Item(457).Insert()

# New row 457: duplicate of the (old) Asterix "1a (guarda lavada)" row, with
# updated date and updated prices.
$ws.Range("A457").Value = 4
$ws.Range("B457").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C457").Value = "Los Lagos"
$ws.Range("D457").Value = 44826
$ws.Range("E457").Value = 10
$ws.Range("F457").Value = 100114001
$ws.Range("G457").Value = "Papa"
$ws.Range("H457").Value = "Asterix"
$ws.Range("I457").Value = "1a (guarda lavada)"
$ws.Range("J457").Value = 150
$ws.Range("K457").Value = 9000
$ws.Range("L457").Value = 9000
$ws.Range("M457").Value = 9000
$ws.Range("N457").Value = "`$/malla 25 kilos"
$ws.Range("O457").Value = "Provincia de Llanquihue"
$ws.Range("P457").Value = 360
$ws.Range("Q457").Value = 25
$ws.Range("R457").Value = "Hortaliza"

# New row 458: duplicate of the (old) Patagonia "1a (guarda)" row, with
# updated date and updated prices.
$ws.Range("A458").Value = 4
$ws.Range("B458").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C458").Value = "Los Lagos"
$ws.Range("D458").Value = 44826
$ws.Range("E458").Value = 10
$ws.Range("F458").Value = 100114001
$ws.Range("G458").Value = "Papa"
$ws.Range("H458").Value = "Patagonia"
$ws.Range("I458").Value = "1a (guarda)"
$ws.Range("J458").Value = 150
$ws.Range("K458").Value = 9000
$ws.Range("L458").Value = 9000
$ws.Range("M458").Value = 9000
$ws.Range("N458").Value = "`$/saco 25 kilos"
$ws.Range("O458").Value = "Provincia de Llanquihue"
$ws.Range("P458").Value = 320
$ws.Range("Q458").Value = 25
$ws.Range("R458").Value = "Hortaliza"

# Match the date number format used by the rest of column D (yyyy-mm-dd-like
# date cells use style index 2 in the original file); copy format from the
# row immediately below so the new date cells render consistently.
$ws.Range("D459").Copy()
$ws.Range("D457:D458").PasteSpecial(-4122)
$excel.CutCopyMode = 0
